$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 4;   I = "%";  J = "Uninterpretable" },
    @{ Row = 26;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 41;  I = "ba"; J = "Appreciation" },
    @{ Row = 47;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 61;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 62;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 66;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 67;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 73;  I = "%";  J = "Uninterpretable" },
    @{ Row = 74;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 75;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 78;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 93;  I = "ba"; J = "Appreciation" },
    @{ Row = 97;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 101; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 114; I = "ba"; J = "Appreciation" },
    @{ Row = 116; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 118; I = "ba"; J = "Appreciation" },
    @{ Row = 126; I = "b";  J = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
